# Update the lattice-multiplication exercise table: every cell keeps its
# layout (problem line, factor-digits line, "----" separator, two lattice
# rows) but gets a fresh pair of factors and matching digit breakdown.
# The table shape (5 rows x 3 columns) is unchanged, so each cell is
# addressed directly and its whole text is replaced in one shot, using a
# vertical-tab (0x0B) as the line-break character so Word re-serializes
# it back to <w:br/> between <w:t> runs, preserving the run formatting
# (sz=32) already on each cell.

$d = $word.ActiveDocument
$vt = [char]0x0B

$d.Tables(1).Cell(1, 1).Range.Text = "41 x 68" + $vt + "  6    8" + $vt + "  ----" + $vt + "4|    |" + $vt + "1|    |"
$d.Tables(1).Cell(1, 2).Range.Text = "73 x 86" + $vt + "  8    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"
$d.Tables(1).Cell(1, 3).Range.Text = "20 x 29" + $vt + "  2    9" + $vt + "  ----" + $vt + "2|    |" + $vt + "0|    |"

$d.Tables(1).Cell(2, 1).Range.Text = "92 x 91" + $vt + "  9    1" + $vt + "  ----" + $vt + "9|    |" + $vt + "2|    |"
$d.Tables(1).Cell(2, 2).Range.Text = "62 x 83" + $vt + "  8    3" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"
$d.Tables(1).Cell(2, 3).Range.Text = "40 x 25" + $vt + "  2    5" + $vt + "  ----" + $vt + "4|    |" + $vt + "0|    |"

$d.Tables(1).Cell(3, 1).Range.Text = "65 x 35" + $vt + "  3    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "5|    |"
$d.Tables(1).Cell(3, 2).Range.Text = "10 x 88" + $vt + "  8    8" + $vt + "  ----" + $vt + "1|    |" + $vt + "0|    |"
$d.Tables(1).Cell(3, 3).Range.Text = "27 x 78" + $vt + "  7    8" + $vt + "  ----" + $vt + "2|    |" + $vt + "7|    |"

$d.Tables(1).Cell(4, 1).Range.Text = "73 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"
$d.Tables(1).Cell(4, 2).Range.Text = "52 x 18" + $vt + "  1    8" + $vt + "  ----" + $vt + "5|    |" + $vt + "2|    |"
$d.Tables(1).Cell(4, 3).Range.Text = "56 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "5|    |" + $vt + "6|    |"

$d.Tables(1).Cell(5, 1).Range.Text = "68 x 89" + $vt + "  8    9" + $vt + "  ----" + $vt + "6|    |" + $vt + "8|    |"
$d.Tables(1).Cell(5, 2).Range.Text = "73 x 38" + $vt + "  3    8" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"
$d.Tables(1).Cell(5, 3).Range.Text = "17 x 70" + $vt + "  7    0" + $vt + "  ----" + $vt + "1|    |" + $vt + "7|    |"
